$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct a false-positive data point for 2014: 148 -> 147, flagged in red.
$ws.Range("B5").Value = 147
$ws.Range("B5").Font.Color = 255

# The formula in B18 (SUM(B1:B17)) recalculates automatically to 2187.

# Move the active selection (as left by the editing session) to E8.
$ws.Range("E8").Select()

# Re-apply the descending sort on column A, now anchored at row 20
# (mirrors the sortState/sortCondition row shift recorded by Excel).
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A20:A36"), 0, 2)
$ws.Sort.SetRange($ws.Range("A20:B36"))
$ws.Sort.Apply()
